{"js": "const replacements = [\n  [\"607\u00d74=2428\", \"839\u00d77=5873\"],\n  [\"519\u00d77=3633\", \"239\u00d72=478\"],\n  [\"287\u00d73=861\", \"759\u00d76=4554\"],\n  [\"291\u00d73=873\", \"332\u00d79=2988\"],\n  [\"261\u00d72=522\", \"772\u00d75=3860\"],\n  [\"180\u00d78=1440\", \"763\u00d78=6104\"],\n  [\"130\u00d78=1040\", \"754\u00d75=3770\"],\n  [\"876\u00d74=3504\", \"343\u00d78=2744\"],\n  [\"982\u00d73=2946\", \"461\u00d77=3227\"],\n  [\"321\u00d76=1926\", \"239\u00d72=478\"],\n  [\"134\u00d79=1206\", \"511\u00d72=1022\"],\n  [\"707\u00d74=2828\", \"402\u00d72=804\"],\n  [\"872\u00d76=5232\", \"577\u00d75=2885\"],\n  [\"523\u00d75=2615\", \"393\u00d72=786\"],\n  [\"694\u00d79=6246\", \"575\u00d73=1725\"],\n  [\"461\u00d73=1383\", \"463\u00d78=3704\"],\n  [\"714\u00d74=2856\", \"380\u00d78=3040\"],\n  [\"911\u00d78=7288\", \"804\u00d73=2412\"],\n  [\"626\u00d78=5008\", \"995\u00d72=1990\"],\n  [\"813\u00d73=2439\", \"557\u00d72=1114\"],\n  [\"394\u00d79=3546\", \"836\u00d76=5016\"],\n  [\"552\u00d78=4416\", \"552\u00d73=1656\"],\n  [\"952\u00d74=3808\", \"592\u00d74=2368\"],\n  [\"312\u00d79=2808\", \"903\u00d73=2709\"],\n  [\"715\u00d76=4290\", \"866\u00d74=3464\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"607\u00d74=2428\", \"839\u00d77=5873\"),\n    @(\"519\u00d77=3633\", \"239\u00d72=478\"),\n    @(\"287\u00d73=861\", \"759\u00d76=4554\"),\n    @(\"291\u00d73=873\", \"332\u00d79=2988\"),\n    @(\"261\u00d72=522\", \"772\u00d75=3860\"),\n    @(\"180\u00d78=1440\", \"763\u00d78=6104\"),\n    @(\"130\u00d78=1040\", \"754\u00d75=3770\"),\n    @(\"876\u00d74=3504\", \"343\u00d78=2744\"),\n    @(\"982\u00d73=2946\", \"461\u00d77=3227\"),\n    @(\"321\u00d76=1926\", \"239\u00d72=478\"),\n    @(\"134\u00d79=1206\", \"511\u00d72=1022\"),\n    @(\"707\u00d74=2828\", \"402\u00d72=804\"),\n    @(\"872\u00d76=5232\", \"577\u00d75=2885\"),\n    @(\"523\u00d75=2615\", \"393\u00d72=786\"),\n    @(\"694\u00d79=6246\", \"575\u00d73=1725\"),\n    @(\"461\u00d73=1383\", \"463\u00d78=3704\"),\n    @(\"714\u00d74=2856\", \"380\u00d78=3040\"),\n    @(\"911\u00d78=7288\", \"804\u00d73=2412\"),\n    @(\"626\u00d78=5008\", \"995\u00d72=1990\"),\n    @(\"813\u00d73=2439\", \"557\u00d72=1114\"),\n    @(\"394\u00d79=3546\", \"836\u00d76=5016\"),\n    @(\"552\u00d78=4416\", \"552\u00d73=1656\"),\n    @(\"952\u00d74=3808\", \"592\u00d74=2368\"),\n    @(\"312\u00d79=2808\", \"903\u00d73=2709\"),\n    @(\"715\u00d76=4290\", \"866\u00d74=3464\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)\n}\n"}
